# Update odds values on Sheet1 to reflect the latest FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (Verona - Monza)
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 3

# Row 11 (Cambuur - Eindhoven FC)
$ws.Range("G11").Value = 1.7
$ws.Range("H11").Value = 4.1
$ws.Range("K11").Value = 2.5
$ws.Range("U11").Value = 1.53
$ws.Range("V11").Value = 2.38
$ws.Range("W11").Value = 10
$ws.Range("AA11").Value = 12
$ws.Range("AB11").Value = 19
$ws.Range("AH11").Value = 17
$ws.Range("AI11").Value = 26
$ws.Range("AO11").Value = 8.5
$ws.Range("AP11").Value = 15
$ws.Range("AS11").Value = 81
$ws.Range("BB11").Value = 126

# Row 25 (Radomiak Radom - Puszcza)
$ws.Range("Q25").Value = 1.93
$ws.Range("R25").Value = 1.93

# Row 35 (Erokspor - Keciorengucu)
$ws.Range("M35").Value = 1.05
$ws.Range("N35").Value = 11

# Row 36 (Rukh Lviv - Karpaty Lviv)
$ws.Range("H36").Value = 3.15
$ws.Range("I36").Value = 3.5
$ws.Range("L36").Value = 4
$ws.Range("P36").Value = 2.82
$ws.Range("U36").Value = 1.75
$ws.Range("V36").Value = 1.85
$ws.Range("AB36").Value = 28
$ws.Range("AD36").Value = 6.1
$ws.Range("AF36").Value = 70
$ws.Range("AI36").Value = 18.5
$ws.Range("AP36").Value = 18
$ws.Range("AQ36").Value = 37
$ws.Range("AR36").Value = 65
$ws.Range("AS36").Value = 200
$ws.Range("AT36").Value = 2.55
$ws.Range("AU36").Value = 6.8
$ws.Range("AW36").Value = 5.4
$ws.Range("AY36").Value = 25
$ws.Range("BA36").Value = 120
$ws.Range("BB36").Value = 300
